$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: add new columns P1 (14) and Q1 (15) ---
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Apply the same style as the existing header cells (bold + border + center/top align) to the new header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows 2-25 ---
# Update changed probability values in columns B,D,E,F,G,H,I,J,K,L,M,N
# and populate the two newly added columns P and Q
# Row 2
$ws.Cells.Item(2, 2).Value = 0.8012425509901675
$ws.Cells.Item(2, 4).Value = 0.03679401558227013
$ws.Cells.Item(2, 5).Value = 0.106911160612746
$ws.Cells.Item(2, 6).Value = 0.7887045552643173
$ws.Cells.Item(2, 7).Value = 0.7050322516239618
$ws.Cells.Item(2, 8).Value = 0.003319578130342005
$ws.Cells.Item(2, 9).Value = 0.01202212491932064
$ws.Cells.Item(2, 10).Value = 0.4762743300367873
$ws.Cells.Item(2, 11).Value = 0.5230135877468847
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 1.751651585294013
$ws.Cells.Item(2, 14).Value = 0.3765695944079113
$ws.Cells.Item(2, 16).Value = 0.9353437802524365
$ws.Cells.Item(2, 17).Value = 0
# Row 3
$ws.Cells.Item(3, 2).Value = 0.7042014500397045
$ws.Cells.Item(3, 4).Value = 0.03439714363844359
$ws.Cells.Item(3, 5).Value = 0.09957054220915662
$ws.Cells.Item(3, 6).Value = 0.7252297804540078
$ws.Cells.Item(3, 7).Value = 0.6484164079478916
$ws.Cells.Item(3, 8).Value = 0.004895276764915746
$ws.Cells.Item(3, 9).Value = 0.01438201941141326
$ws.Cells.Item(3, 10).Value = 0.4551740425176547
$ws.Cells.Item(3, 11).Value = 0.5174525311049116
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 1.524070879469235
$ws.Cells.Item(3, 14).Value = 0.328501578772503
$ws.Cells.Item(3, 16).Value = 0.9674108242651727
$ws.Cells.Item(3, 17).Value = 0
# Row 4
$ws.Cells.Item(4, 2).Value = 0.6440445907797425
$ws.Cells.Item(4, 4).Value = 0.03292061904922505
$ws.Cells.Item(4, 5).Value = 0.09506854749377491
$ws.Cells.Item(4, 6).Value = 0.6870930510959425
$ws.Cells.Item(4, 7).Value = 0.6144162392875643
$ws.Cells.Item(4, 8).Value = 0.006056907437054726
$ws.Cells.Item(4, 9).Value = 0.01602889380145722
$ws.Cells.Item(4, 10).Value = 0.4426996745214922
$ws.Cells.Item(4, 11).Value = 0.5140615073866108
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 1.384964161304907
$ws.Cells.Item(4, 14).Value = 0.2991281207693675
$ws.Cells.Item(4, 16).Value = 0.9876052885544873
$ws.Cells.Item(4, 17).Value = 0
# Row 5
$ws.Cells.Item(5, 2).Value = 0.6180696972154465
$ws.Cells.Item(5, 4).Value = 0.03234469940153062
$ws.Cells.Item(5, 5).Value = 0.09319500910383383
$ws.Cells.Item(5, 6).Value = 0.6712050266772209
$ws.Cells.Item(5, 7).Value = 0.6001389723448938
$ws.Cells.Item(5, 8).Value = 0.006581088058957241
$ws.Cells.Item(5, 9).Value = 0.0168405194403376
$ws.Cells.Item(5, 10).Value = 0.4374228874667949
$ws.Cells.Item(5, 11).Value = 0.5120541144153883
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 1.329073243532378
$ws.Cells.Item(5, 14).Value = 0.2874428964219504
$ws.Cells.Item(5, 16).Value = 0.9956215213043986
$ws.Cells.Item(5, 17).Value = 0
# Row 6
$ws.Cells.Item(6, 2).Value = 0.6121373076084353
$ws.Cells.Item(6, 4).Value = 0.03228233030564809
$ws.Cells.Item(6, 5).Value = 0.09283949923132973
$ws.Cells.Item(6, 6).Value = 0.6679157386825949
$ws.Cells.Item(6, 7).Value = 0.5970470102357694
$ws.Cells.Item(6, 8).Value = 0.006675750232228772
$ws.Cells.Item(6, 9).Value = 0.01709815118358815
$ws.Cells.Item(6, 10).Value = 0.4361773877291029
$ws.Cells.Item(6, 11).Value = 0.5109512150174886
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 1.320615975895379
$ws.Cells.Item(6, 14).Value = 0.2858152245657095
$ws.Cells.Item(6, 16).Value = 0.996550900958022
$ws.Cells.Item(6, 17).Value = 0
# Row 7
$ws.Cells.Item(7, 2).Value = 0.6392809448046535
$ws.Cells.Item(7, 4).Value = 0.03300452865568104
$ws.Cells.Item(7, 5).Value = 0.09492307253267085
$ws.Cells.Item(7, 6).Value = 0.6850683954789929
$ws.Cells.Item(7, 7).Value = 0.612224665495134
$ws.Cells.Item(7, 8).Value = 0.00607609552612326
$ws.Cells.Item(7, 9).Value = 0.0163501821969172
$ws.Cells.Item(7, 10).Value = 0.4416020238621314
$ws.Cells.Item(7, 11).Value = 0.5119331701433012
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 1.386444360102274
$ws.Cells.Item(7, 14).Value = 0.2998214046675827
$ws.Cells.Item(7, 16).Value = 0.9866120364630664
$ws.Cells.Item(7, 17).Value = 0
# Row 8
$ws.Cells.Item(8, 2).Value = 0.762009481381142
$ws.Cells.Item(8, 4).Value = 0.03609250277994036
$ws.Cells.Item(8, 5).Value = 0.104222162626175
$ws.Cells.Item(8, 6).Value = 0.7642391989776627
$ws.Cells.Item(8, 7).Value = 0.6826980110350149
$ws.Cells.Item(8, 8).Value = 0.003833527875953369
$ws.Cells.Item(8, 9).Value = 0.01316720806154592
$ws.Cells.Item(8, 10).Value = 0.4675412864032467
$ws.Cells.Item(8, 11).Value = 0.5183249170279041
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 1.676011721128219
$ws.Cells.Item(8, 14).Value = 0.3610987306863933
$ws.Cells.Item(8, 16).Value = 0.944915823748584
$ws.Cells.Item(8, 17).Value = 0
# Row 9
$ws.Cells.Item(9, 2).Value = 1.006093855367737
$ws.Cells.Item(9, 4).Value = 0.04194967351335777
$ws.Cells.Item(9, 5).Value = 0.1228811740709015
$ws.Cells.Item(9, 6).Value = 0.9301471477332086
$ws.Cells.Item(9, 7).Value = 0.8312735601493131
$ws.Cells.Item(9, 8).Value = 0.001055506176652177
$ws.Cells.Item(9, 9).Value = 0.008091898718104673
$ws.Cells.Item(9, 10).Value = 0.5243936043180355
$ws.Cells.Item(9, 11).Value = 0.5348142255468638
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 2.246186458399364
$ws.Cells.Item(9, 14).Value = 0.4810526220424833
$ws.Cells.Item(9, 16).Value = 0.8691473619364878
$ws.Cells.Item(9, 17).Value = 0
# Row 10
$ws.Cells.Item(10, 2).Value = 1.163527868966042
$ws.Cells.Item(10, 4).Value = 0.04688618637849373
$ws.Cells.Item(10, 5).Value = 0.1362968900760728
$ws.Cells.Item(10, 6).Value = 1.046546106748167
$ws.Cells.Item(10, 7).Value = 0.9333836320166711
$ws.Cells.Item(10, 8).Value = 0.0003685001762820939
$ws.Cells.Item(10, 9).Value = 0.005692458475423656
$ws.Cells.Item(10, 10).Value = 0.5630063660035205
$ws.Cells.Item(10, 11).Value = 0.5375814601542714
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 2.675732665693232
$ws.Cells.Item(10, 14).Value = 0.5509390373472343
$ws.Cells.Item(10, 16).Value = 0.8121023760836295
$ws.Cells.Item(10, 17).Value = 0
# Row 11
$ws.Cells.Item(11, 2).Value = 1.083439601052703
$ws.Cells.Item(11, 4).Value = 0.05490395559630912
$ws.Cells.Item(11, 5).Value = 0.1483445402074217
$ws.Cells.Item(11, 6).Value = 1.009343489586385
$ws.Cells.Item(11, 7).Value = 0.8807980204295234
$ws.Cells.Item(11, 8).Value = 0.01896810401230553
$ws.Cells.Item(11, 9).Value = 0.00589783363624985
$ws.Cells.Item(11, 10).Value = 0.5307502274132219
$ws.Cells.Item(11, 11).Value = 0.4669490356675148
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = 2.908502787659074
$ws.Cells.Item(11, 14).Value = 0.4236816437197319
$ws.Cells.Item(11, 16).Value = 0.7625350508016357
$ws.Cells.Item(11, 17).Value = 0
# Row 12
$ws.Cells.Item(12, 2).Value = 0.9924247390585776
$ws.Cells.Item(12, 4).Value = 0.06108222542267328
$ws.Cells.Item(12, 5).Value = 0.1705254023006759
$ws.Cells.Item(12, 6).Value = 0.9559498712806658
$ws.Cells.Item(12, 7).Value = 0.8185999071603618
$ws.Cells.Item(12, 8).Value = 0.05772099078665605
$ws.Cells.Item(12, 9).Value = 0.005974074574037758
$ws.Cells.Item(12, 10).Value = 0.4970809309256907
$ws.Cells.Item(12, 11).Value = 0.4159036992312366
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = 3.011422913076103
$ws.Cells.Item(12, 14).Value = 0.3155406177074838
$ws.Cells.Item(12, 16).Value = 0.7445386040455677
$ws.Cells.Item(12, 17).Value = 0
# Row 13
$ws.Cells.Item(13, 2).Value = 0.8814266577294916
$ws.Cells.Item(13, 4).Value = 0.06626905710572117
$ws.Cells.Item(13, 5).Value = 0.2009187176667311
$ws.Cells.Item(13, 6).Value = 0.8838918602149164
$ws.Cells.Item(13, 7).Value = 0.7420584094934242
$ws.Cells.Item(13, 8).Value = 0.1136153535023965
$ws.Cells.Item(13, 9).Value = 0.006395443611055818
$ws.Cells.Item(13, 10).Value = 0.4586274430448611
$ws.Cells.Item(13, 11).Value = 0.373623637278051
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 3.022532871309807
$ws.Cells.Item(13, 14).Value = 0.2178062313814877
$ws.Cells.Item(13, 16).Value = 0.7471521858916255
$ws.Cells.Item(13, 17).Value = 0
# Row 14
$ws.Cells.Item(14, 2).Value = 0.7946768281954633
$ws.Cells.Item(14, 4).Value = 0.06951613411839475
$ws.Cells.Item(14, 5).Value = 0.2276398397075781
$ws.Cells.Item(14, 6).Value = 0.8249853919861749
$ws.Cells.Item(14, 7).Value = 0.6820706896780138
$ws.Cells.Item(14, 8).Value = 0.1631016112915376
$ws.Cells.Item(14, 9).Value = 0.00695586281751126
$ws.Cells.Item(14, 10).Value = 0.4295405300617148
$ws.Cells.Item(14, 11).Value = 0.3485368583523538
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 2.989433745687165
$ws.Cells.Item(14, 14).Value = 0.1569662775820788
$ws.Cells.Item(14, 16).Value = 0.7592815465547886
$ws.Cells.Item(14, 17).Value = 0
# Row 15
$ws.Cells.Item(15, 2).Value = 0.767720197037022
$ws.Cells.Item(15, 4).Value = 0.07005762252595815
$ws.Cells.Item(15, 5).Value = 0.2341334468444742
$ws.Cells.Item(15, 6).Value = 0.8061470824297459
$ws.Cells.Item(15, 7).Value = 0.6637484649469059
$ws.Cells.Item(15, 8).Value = 0.1756766501975591
$ws.Cells.Item(15, 9).Value = 0.007316316900131348
$ws.Cells.Item(15, 10).Value = 0.4211178614518332
$ws.Cells.Item(15, 11).Value = 0.343130444896051
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 2.962558568609751
$ws.Cells.Item(15, 14).Value = 0.1426926381376816
$ws.Cells.Item(15, 16).Value = 0.7652856837013182
$ws.Cells.Item(15, 17).Value = 0
# Row 16
$ws.Cells.Item(16, 2).Value = 0.7235463410428054
$ws.Cells.Item(16, 4).Value = 0.06671757788831911
$ws.Cells.Item(16, 5).Value = 0.2209593725202907
$ws.Cells.Item(16, 6).Value = 0.7696877330882756
$ws.Cells.Item(16, 7).Value = 0.634458004817958
$ws.Cells.Item(16, 8).Value = 0.1632912075477719
$ws.Cells.Item(16, 9).Value = 0.008456716015698085
$ws.Cells.Item(16, 10).Value = 0.4116637167082615
$ws.Cells.Item(16, 11).Value = 0.3497346721846597
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = 2.774188289049619
$ws.Cells.Item(16, 14).Value = 0.1359058518869318
$ws.Cells.Item(16, 16).Value = 0.7859829670460705
$ws.Cells.Item(16, 17).Value = 0
# Row 17
$ws.Cells.Item(17, 2).Value = 0.7346595126556679
$ws.Cells.Item(17, 4).Value = 0.06243044397895403
$ws.Cells.Item(17, 5).Value = 0.1955876313923994
$ws.Cells.Item(17, 6).Value = 0.7730427114577196
$ws.Cells.Item(17, 7).Value = 0.6437565783588894
$ws.Cells.Item(17, 8).Value = 0.1257161231439028
$ws.Cells.Item(17, 9).Value = 0.009062088133203083
$ws.Cells.Item(17, 10).Value = 0.4198655046197644
$ws.Cells.Item(17, 11).Value = 0.3674823733303469
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = 2.648058658542283
$ws.Cells.Item(17, 14).Value = 0.1611714339297237
$ws.Cells.Item(17, 16).Value = 0.795105036727012
$ws.Cells.Item(17, 17).Value = 0
# Row 18
$ws.Cells.Item(18, 2).Value = 0.7996902381757138
$ws.Cells.Item(18, 4).Value = 0.05698619009103112
$ws.Cells.Item(18, 5).Value = 0.1632995921046998
$ws.Cells.Item(18, 6).Value = 0.8124080252069774
$ws.Cells.Item(18, 7).Value = 0.6892508425918891
$ws.Cells.Item(18, 8).Value = 0.07288972416720441
$ws.Cells.Item(18, 9).Value = 0.008864126068223399
$ws.Cells.Item(18, 10).Value = 0.4450917028080426
$ws.Cells.Item(18, 11).Value = 0.4010131775237049
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 2.558848452549057
$ws.Cells.Item(18, 14).Value = 0.224397019219154
$ws.Cells.Item(18, 16).Value = 0.8001660537482032
$ws.Cells.Item(18, 17).Value = 0
# Row 19
$ws.Cells.Item(19, 2).Value = 0.8995077354829846
$ws.Cells.Item(19, 4).Value = 0.05160104066895599
$ws.Cells.Item(19, 5).Value = 0.1386954244040783
$ws.Cells.Item(19, 6).Value = 0.8756131361006538
$ws.Cells.Item(19, 7).Value = 0.7587482806641646
$ws.Cells.Item(19, 8).Value = 0.02738982231841902
$ws.Cells.Item(19, 9).Value = 0.008599173557838569
$ws.Cells.Item(19, 10).Value = 0.4810944385808114
$ws.Cells.Item(19, 11).Value = 0.4470472425313048
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 2.510534313824365
$ws.Cells.Item(19, 14).Value = 0.3270931465489753
$ws.Cells.Item(19, 16).Value = 0.8088299811026642
$ws.Cells.Item(19, 17).Value = 0
# Row 20
$ws.Cells.Item(20, 2).Value = 1.10768617238196
$ws.Cells.Item(20, 4).Value = 0.04592235938982725
$ws.Cells.Item(20, 5).Value = 0.1324122509631385
$ws.Cells.Item(20, 6).Value = 1.009491881651073
$ws.Cells.Item(20, 7).Value = 0.8994886846741679
$ws.Cells.Item(20, 8).Value = 0.0004657635245446379
$ws.Cells.Item(20, 9).Value = 0.007156889192149407
$ws.Cells.Item(20, 10).Value = 0.549245841162687
$ws.Cells.Item(20, 11).Value = 0.5301437639567403
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = 2.56972523160286
$ws.Cells.Item(20, 14).Value = 0.5343848369068525
$ws.Cells.Item(20, 16).Value = 0.8239632134093782
$ws.Cells.Item(20, 17).Value = 0
# Row 21
$ws.Cells.Item(21, 2).Value = 1.256531023579072
$ws.Cells.Item(21, 4).Value = 0.04873853938809702
$ws.Cells.Item(21, 5).Value = 0.1439444277369
$ws.Cells.Item(21, 6).Value = 1.118044553005234
$ws.Cells.Item(21, 7).Value = 0.9986982088190928
$ws.Cells.Item(21, 8).Value = [double]"5.272392679778193E-06"
$ws.Cells.Item(21, 9).Value = 0.00535880293281199
$ws.Cells.Item(21, 10).Value = 0.5894819939099563
$ws.Cells.Item(21, 11).Value = 0.5468681737749606
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = 2.894824252638216
$ws.Cells.Item(21, 14).Value = 0.6185145709790874
$ws.Cells.Item(21, 16).Value = 0.7875676946088657
$ws.Cells.Item(21, 17).Value = 0
# Row 22
$ws.Cells.Item(22, 2).Value = 1.35065857941072
$ws.Cells.Item(22, 4).Value = 0.05075310981062842
$ws.Cells.Item(22, 5).Value = 0.1514720352832506
$ws.Cells.Item(22, 6).Value = 1.187540369396785
$ws.Cells.Item(22, 7).Value = 1.061712106181531
$ws.Cells.Item(22, 8).Value = [double]"5.747215913776138E-05"
$ws.Cells.Item(22, 9).Value = 0.004106483773627545
$ws.Cells.Item(22, 10).Value = 0.6149207369475818
$ws.Cells.Item(22, 11).Value = 0.5561933766067
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 3.107933624933338
$ws.Cells.Item(22, 14).Value = 0.6625816636719151
$ws.Cells.Item(22, 16).Value = 0.7635965843381847
$ws.Cells.Item(22, 17).Value = 0
# Row 23
$ws.Cells.Item(23, 2).Value = 1.3058145279727
$ws.Cells.Item(23, 4).Value = 0.0495596941936256
$ws.Cells.Item(23, 5).Value = 0.1475576349648953
$ws.Cells.Item(23, 6).Value = 1.152445183696813
$ws.Cells.Item(23, 7).Value = 1.030304126099949
$ws.Cells.Item(23, 8).Value = [double]"5.755710401844638E-06"
$ws.Cells.Item(23, 9).Value = 0.004413176288946552
$ws.Cells.Item(23, 10).Value = 0.6024675629493146
$ws.Cells.Item(23, 11).Value = 0.5536194157049792
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = 2.991316465077944
$ws.Cells.Item(23, 14).Value = 0.6379983998644434
$ws.Cells.Item(23, 16).Value = 0.7773482352763939
$ws.Cells.Item(23, 17).Value = 0
# Row 24
$ws.Cells.Item(24, 2).Value = 1.12667774243036
$ws.Cells.Item(24, 4).Value = 0.04523005649470235
$ws.Cells.Item(24, 5).Value = 0.1328825781046106
$ws.Cells.Item(24, 6).Value = 1.019129684364103
$ws.Cells.Item(24, 7).Value = 0.9104470362245394
$ws.Cells.Item(24, 8).Value = 0.0003145893420370971
$ws.Cells.Item(24, 9).Value = 0.006562043901560344
$ws.Cells.Item(24, 10).Value = 0.5549360197247921
$ws.Cells.Item(24, 11).Value = 0.5402672022968318
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 2.55747288274074
$ws.Cells.Item(24, 14).Value = 0.5471215785739645
$ws.Cells.Item(24, 16).Value = 0.8290758996807437
$ws.Cells.Item(24, 17).Value = 0
# Row 25
$ws.Cells.Item(25, 2).Value = 0.9324376668461696
$ws.Cells.Item(25, 4).Value = 0.0405414038245091
$ws.Cells.Item(25, 5).Value = 0.1175853753700329
$ws.Cells.Item(25, 6).Value = 0.8810833403112071
$ws.Cells.Item(25, 7).Value = 0.786627875717727
$ws.Cells.Item(25, 8).Value = 0.001628145225335831
$ws.Cells.Item(25, 9).Value = 0.00979175951476563
$ws.Cells.Item(25, 10).Value = 0.506677216542343
$ws.Cells.Item(25, 11).Value = 0.5265842023972347
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 2.09515437728524
$ws.Cells.Item(25, 14).Value = 0.4499944262765183
$ws.Cells.Item(25, 16).Value = 0.887438859694484
$ws.Cells.Item(25, 17).Value = 0
